# Generate Report for Handoff
# Adds a new handoff row (29659e31-4389-459a-b885-d57447a22a96) to the
# Overview sheet and to each locale sheet (zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$fileId   = "29659e31-4389-459a-b885-d57447a22a96"
$fileHash = "eb4f003860cfb5f96f717968c750eb4cb2b4d2a2"

$mdName      = "$fileId.md"
$zhXlfName   = "$fileId.$fileHash.zh-cn.xlf"
$deXlfName   = "$fileId.$fileHash.de-de.xlf"

$mdUrl    = "https://github.com/OpenLocalizationTest/oltest/blob/86a2f557c32963927a6ce6f3671a2b72086c7e4b/e2e/$mdName"
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5ce216930cc8e46fb6efa299885909994cdbdbc4/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$zhXlfName"
$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/661591660b77433ad0c3c8bbadf363b8e8cfa3c8/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$deXlfName"

$readyStatus = "Ready for handoff"
$emptyDate   = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------------
# Overview sheet: File Name | zh-cn | de-de | Latest Handoff Date
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

$overview.Range("B3").Value = $readyStatus
$overview.Range("C3").Value = $readyStatus
$overview.Range("D3").Value = "2016-03-21 16:38:30"
$overview.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$overview.Hyperlinks.Add($overview.Range("A3"), $mdUrl, $null, $null, $mdName) | Out-Null

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("B3").Value = ".md"
$zhcn.Range("C3").Value = $readyStatus
$zhcn.Range("E3").Value = "2016-03-21 16:38:26"
$zhcn.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zhcn.Range("H3").Value = $emptyDate
$zhcn.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zhcn.Range("J3").Value = "Include"

$zhcn.Hyperlinks.Add($zhcn.Range("A3"), $mdUrl, $null, $null, $mdName) | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("D3"), $zhXlfUrl, $null, $null, $zhXlfName) | Out-Null

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("B3").Value = ".md"
$dede.Range("C3").Value = $readyStatus
$dede.Range("E3").Value = "2016-03-21 16:38:30"
$dede.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$dede.Range("H3").Value = $emptyDate
$dede.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$dede.Range("J3").Value = "Include"

$dede.Hyperlinks.Add($dede.Range("A3"), $mdUrl, $null, $null, $mdName) | Out-Null
$dede.Hyperlinks.Add($dede.Range("D3"), $deXlfUrl, $null, $null, $deXlfName) | Out-Null

Write-Host "Handoff row added for $fileId"
